# Appends 10 new translation rows (Sl. No. 46-55) to the "translations"
# sheet, continuing the existing English<->Santhali table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46, "English → Santhali", "adada", "ᱴᱟᱹᱯᱩ ᱨᱮ", ""),
    @(47, "English → Santhali", "asas", "ᱵᱮᱥ", ""),
    @(48, "English → Santhali", "asasas", "ᱟᱥᱟᱥᱟᱥ", ""),
    @(49, "English → Santhali", "adadadad", "ᱟᱰᱟᱰᱟᱰᱟᱰᱟ", ""),
    @(50, "English → Santhali", "adadad", "ᱵᱟᱵᱟ", "ᱵᱟᱵᱟᱢᱢᱨ"),
    @(51, "English → Santhali", "Prerna Pagal Hai", "ᱯᱨᱮᱨᱱᱟ ᱫᱚ ᱯᱟᱜᱽᱞᱟ ᱜᱮᱭᱟᱭ", ""),
    @(52, "English → Santhali", "Prerna is a good girl", "ᱯᱨᱮᱨᱱᱟ ᱫᱚ ᱱᱟᱯᱟᱭ ᱠᱩᱲᱤ ᱠᱟᱱᱟᱭ", "ᱯᱨᱮᱨᱱᱟ ᱫᱚ ᱱᱟᱯᱟᱭ ᱠᱩᱲᱤ ᱠᱟᱱᱟᱭᱢᱣᱤ"),
    @(53, "English → Santhali", "Hi, My name is Nandu aka Nandani", "ᱦᱟᱭ, ᱤᱧᱟᱜ ᱧᱩᱛᱩᱢ ᱫᱚ ᱱᱟᱱᱰᱩ ᱟᱨᱠᱟ ᱱᱟᱱᱰᱟᱱᱤ", "ᱦᱟᱭ, ᱤᱧᱟᱜ ᱧᱩᱛᱩᱢ ᱫᱚ ᱱᱟᱱᱰᱩ ᱟᱨᱠᱟ ᱱᱟᱱᱰᱟᱱᱤᱠᱡᱢ"),
    @(54, "English → Santhali", "adad", "ᱵᱟᱵᱟ", ""),
    @(55, "Santhali → English", "ᱣᱡᱧᱧ", "Wjng", "")
)

$startRow = 47
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    # Column E ("Updated Text") is blank for most of these rows. Assigning
    # an empty string clears/omits the cell (same as leaving it untouched),
    # which matches how the workbook treats blank cells in this column.
    if ($row[4] -ne "") {
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
}
